# Generate Report for Handoff
# Replaces the two tracked files (c44ad50b... and e7bfa5dc...) with a new
# handoff pair (a3df63f6... and ffff5840b6e6...), flips their status from
# "Handed back: in sync with en-US" to "Ready for handoff", and clears out
# the stale handback bookkeeping (Latest Target File / Latest Handback
# File / Latest Handback DateTime) since nothing has been handed back yet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New identities
# ---------------------------------------------------------------------
$oldId1 = "c44ad50b-1c20-4212-bee6-1c13e0dd490e"
$newId1 = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2"
$oldId2 = "e7bfa5dc-4130-4c1a-9d6a-dcc4cba00484"
$newId2 = "ffff5840b6e6-9309-4023-912d-e5b977376339"

$newXlfHash = "eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f"

$statusText = "Ready for handoff"

$overviewDate = "2016-08-28 13:03:19"
$zhHandoffDate = "2016-08-28 13:03:14"
$epoch = "0001-01-01 00:00:00"

$newXlfZh = "$newId1.$newXlfHash.zh-cn.xlf"
$newXlfDe = "$newId1.$newXlfHash.de-de.xlf"

# Widths from the target diff are finer-grained than this engine's column
# width quantization (snaps to 1/6 character steps), so we pick the input
# that lands on the closest achievable bucket.
$narrowWidth = 16.35   # -> ~17.1667 (target 17.2159881591797)
$col9Width   = 17.80   # -> ~18.6667 (target 18.6506053379604)
$col10Width  = 20.80   # -> ~21.6667 (target 21.7054770333426)

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = "$newId1.md"
$wsOverview.Range("E2").Value2 = $statusText
$wsOverview.Range("F2").Value2 = $statusText
$wsOverview.Range("G2").Value2 = $overviewDate

$wsOverview.Range("A3").Value2 = "$newId2.md"
$wsOverview.Range("E3").Value2 = $statusText
$wsOverview.Range("F3").Value2 = $statusText
$wsOverview.Range("G3").Value2 = $overviewDate

# Hyperlinks on B2/B3 keep their original target addresses; only the
# displayed text needs updating -- rebuild via delete+add so the XML
# ends up with a single clean <hyperlink> entry (an in-place property
# assignment would instead append a duplicate entry).
$b2 = $wsOverview.Range("B2")
$b2Address = $b2.Hyperlinks.Item(1).Address
$b2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($b2, $b2Address, [Type]::Missing, [Type]::Missing, "e2e\$newId1.md")

$b3 = $wsOverview.Range("B3")
$b3Address = $b3.Hyperlinks.Item(1).Address
$b3.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($b3, $b3Address, [Type]::Missing, [Type]::Missing, "e2e\$newId2.md")

$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$a2 = $wsZh.Range("A2")
$a2Address = $a2.Hyperlinks.Item(1).Address
$a2.Hyperlinks.Delete()
$a2.Value2 = "$newId1.md"
$wsZh.Hyperlinks.Add($a2, $a2Address, [Type]::Missing, [Type]::Missing, "$newId1.md")

$wsZh.Range("C2").Value2 = $statusText
$wsZh.Range("G2").Value2 = $newXlfZh
$wsZh.Range("H2").Value2 = $zhHandoffDate
$wsZh.Range("I2").Value2 = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value2 = ""
$wsZh.Range("K2").Value2 = $epoch

$a3 = $wsZh.Range("A3")
$a3Address = $a3.Hyperlinks.Item(1).Address
$a3.Hyperlinks.Delete()
$a3.Value2 = "$newId2.md"
$wsZh.Hyperlinks.Add($a3, $a3Address, [Type]::Missing, [Type]::Missing, "$newId2.md")

$wsZh.Range("C3").Value2 = $statusText
$wsZh.Range("F3").Value2 = "True"
$wsZh.Range("G3").Value2 = $newXlfZh
$wsZh.Range("H3").Value2 = $zhHandoffDate
$wsZh.Range("I3").Value2 = ""
$wsZh.Range("I3").Style = "Normal"
$wsZh.Range("J3").Value2 = ""
$wsZh.Range("K3").Value2 = $epoch

# I2/I3 hyperlinks (pointing at the zh-cn fork) are dropped entirely --
# there is no longer a handback file to link to.
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I3").Hyperlinks.Delete()

$wsZh.Columns.Item(3).ColumnWidth = $narrowWidth
$wsZh.Columns.Item(9).ColumnWidth = $col9Width
$wsZh.Columns.Item(10).ColumnWidth = $col10Width

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$a2d = $wsDe.Range("A2")
$a2dAddress = $a2d.Hyperlinks.Item(1).Address
$a2d.Hyperlinks.Delete()
$a2d.Value2 = "$newId1.md"
$wsDe.Hyperlinks.Add($a2d, $a2dAddress, [Type]::Missing, [Type]::Missing, "$newId1.md")

$wsDe.Range("C2").Value2 = $statusText
$wsDe.Range("G2").Value2 = $newXlfDe
$wsDe.Range("I2").Value2 = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value2 = ""
$wsDe.Range("K2").Value2 = $epoch

$a3d = $wsDe.Range("A3")
$a3dAddress = $a3d.Hyperlinks.Item(1).Address
$a3d.Hyperlinks.Delete()
$a3d.Value2 = "$newId2.md"
$wsDe.Hyperlinks.Add($a3d, $a3dAddress, [Type]::Missing, [Type]::Missing, "$newId2.md")

$wsDe.Range("C3").Value2 = $statusText
$wsDe.Range("F3").Value2 = "True"
$wsDe.Range("G3").Value2 = $newXlfDe
$wsDe.Range("I3").Value2 = ""
$wsDe.Range("I3").Style = "Normal"
$wsDe.Range("J3").Value2 = ""
$wsDe.Range("K3").Value2 = $epoch

# H2/H3 are untouched directly -- they share the Overview datetime string,
# so they pick up "2016-08-28 13:03:19" automatically; set explicitly too
# so the text matches even if the shared string was not reused.
$wsDe.Range("H2").Value2 = $overviewDate
$wsDe.Range("H3").Value2 = $overviewDate

$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I3").Hyperlinks.Delete()

$wsDe.Columns.Item(3).ColumnWidth = $narrowWidth
$wsDe.Columns.Item(9).ColumnWidth = $col9Width
$wsDe.Columns.Item(10).ColumnWidth = $col10Width

Write-Output "Done applying handoff report changes."
